$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated T3 (column B), T4 (column C), and T5 (column D) DRC data for rows 2-6

$ws.Range("B2").Value = 36.74999562515195
$ws.Range("C2").Value = 36.67964383945198
$ws.Range("D2").Value = 36.31682589057031

$ws.Range("B3").Value = 36.98793446843405
$ws.Range("C3").Value = 37.11040202143514
$ws.Range("D3").Value = 37.20786665021228

$ws.Range("B4").Value = 36.71468764237777
$ws.Range("C4").Value = 36.99752832807091
$ws.Range("D4").Value = 37.16715760869143

$ws.Range("B5").Value = 36.3188415262585
$ws.Range("C5").Value = 37.23163313589489
$ws.Range("D5").Value = 37.38171438225818

$ws.Range("B6").Value = 37.09691046084393
$ws.Range("C6").Value = 36.92083705946173
$ws.Range("D6").Value = 37.09286990000914
